$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q4" sheet right after "总计" (before old "2022-Q3"),
#    by duplicating the "2022-Q3" sheet so header styling/borders match the
#    existing quarterly sheets exactly.
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)
$q4Sheet = $q3Sheet
$q4Sheet.Name = "2022-Q4"

$q4Data = @(
    @("0", "002450", "平安睿享文娱灵活配置混合A", "3.85", "94.03", "4.38", "0.1686", "4"),
    @("1", "010994", "博时创新经济混合A", "3.80", "86.65", "3.66", "0.1391", "8"),
    @("2", "519664", "银河美丽优萃混合A", "2.47", "93.91", "5.23", "0.1292", "7"),
    @("3", "519651", "银河转型增长主题灵活配置混合", "4.25", "84.25", "2.69", "0.1143", "10"),
    @("4", "002451", "平安睿享文娱灵活配置混合C", "1.97", "94.03", "4.38", "0.0863", "4"),
    @("5", "010995", "博时创新经济混合C", "1.71", "86.65", "3.66", "0.0626", "8"),
    @("6", "016336", "博时卓远成长一年持有期股票A", "1.63", "51.15", "3.43", "0.0559", "8"),
    @("7", "013687", "平安成长龙头1年持有混合A", "1.21", "94.87", "4.11", "0.0497", "7"),
    @("8", "519665", "银河美丽优萃混合C", "0.72", "93.91", "5.23", "0.0377", "7"),
    @("9", "013688", "平安成长龙头1年持有混合C", "0.51", "94.87", "4.11", "0.0210", "7"),
    @("10", "016337", "博时卓远成长一年持有期股票C", "0.46", "51.15", "3.43", "0.0158", "8"),
    @("11", "350007", "天治趋势精选混合", "0.39", "93.83", "3.36", "0.0131", "7"),
    @("12", "003105", "光大保德信永鑫灵活配置混合A", "0.93", "24.66", "0.92", "0.0086", "10"),
    @("13", "001464", "光大保德信鼎鑫灵活配置混合A", "0.89", "21.94", "0.82", "0.0073", "10"),
    @("14", "001823", "光大保德信鼎鑫灵活配置混合C", "0.25", "21.94", "0.82", "0.0020", "10"),
    @("15", "003106", "光大保德信永鑫灵活配置混合C", "0.09", "24.66", "0.92", "0.0008", "10")
)

# Extend the copied sheet (currently 2 data rows) down to 16 data rows,
# cloning the style of row 2 (index col A, plain data cols B:H) for the
# newly-needed rows.
for ($r = 4; $r -le 17; $r++) {
    $q4Sheet.Range("A2").Copy($q4Sheet.Range("A$r"))
    $q4Sheet.Range("B2:H2").Copy($q4Sheet.Range("B$r`:H$r"))
}

# Columns B:G hold text (fund code/name/size/position numbers stored as
# strings in the source data); force text formatting before writing so
# leading zeros and decimal-looking strings are preserved verbatim.
$q4Sheet.Range("B2:G17").NumberFormat = "@"

$r = 2
foreach ($row in $q4Data) {
    $q4Sheet.Cells.Item($r, 1).Value = [int]$row[0]
    $q4Sheet.Cells.Item($r, 2).Value = $row[1]
    $q4Sheet.Cells.Item($r, 3).Value = $row[2]
    $q4Sheet.Cells.Item($r, 4).Value = $row[3]
    $q4Sheet.Cells.Item($r, 5).Value = $row[4]
    $q4Sheet.Cells.Item($r, 6).Value = $row[5]
    $q4Sheet.Cells.Item($r, 7).Value = $row[6]
    $q4Sheet.Cells.Item($r, 8).Value = [int]$row[7]
    $r = $r + 1
}

# The temporary "@" text format is no longer needed now that the strings
# are committed as text; drop it so these data cells go back to the
# unstyled look the source rows use (matches the other quarter sheets).
$q4Sheet.Range("B2:G17").ClearFormats()

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a row for 2022-Q4 above the
#    existing 2022-Q3 row (pushing 2022-Q3 / 2022-Q2 down one row each).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(3).Insert()

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 16
$totalSheet.Range("D2").Value = 0.91

# Row-insert doesn't fully clone the index column's border style, so
# re-stamp A3 from A2 (both use the same bold+bordered "index" style)
# before writing its value.
$totalSheet.Range("A2").Copy($totalSheet.Range("A3"))
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.01

$totalSheet.Range("A4").Value = 2

Write-Output "edit complete"
